# "page model correction ratio"
# Insert a new "page/record" column (between "model type" and "bidirectional"),
# and insert a new data row describing the "record"-level metrics for the
# first_try / LSTM run (row 2 becomes the "page" metrics, the new row 3
# becomes the "record" metrics).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- structural edits -------------------------------------------------
# Insert a new column before D (old D "bidirectional" -> E, etc.)
$ws.Range("D1").EntireColumn.Insert()

# Insert a new row before row 3 (old rows 3-7 -> 4-8)
$ws.Range("A3").EntireRow.Insert()

# --- header row ---------------------------------------------------------
$ws.Range("D1").Value = "page/record"

# --- row 2: existing first_try/LSTM run is now labelled "page" ----------
$ws.Range("D2").Value = "page"

# --- row 3 (new): matching "record" level metrics for the same run ------
$ws.Range("B3").Value = "first_try"
$ws.Range("C3").Value = "LSTM"
$ws.Range("D3").Value = "record"
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = 0.05
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 20
$ws.Range("J3").Value = 0.079865917559999999
$ws.Range("L3").Value = 0.89944000000000002
$ws.Range("M3").Value = 0.89479797271749995

# --- column D width (new column, not a bestFit column) ------------------
$ws.Range("D1").ColumnWidth = 10.17

# --- restore selection shown in the saved workbook -----------------------
[void]$ws.Range("M4").Select()
